$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row above row 18, shifting existing rows 18-36 down to 19-37
$ws.Rows.Item(18).Insert()

# The newly inserted row is blank/unformatted; pull the borders/styles from the
# row pushed down to 19 (the old blank spacer row 18, which already carried the
# correct per-column styles s=7,8,8,8,9,9,10,11) so the new row matches it.
$ws.Range("A19:H19").Copy() | Out-Null
$ws.Range("A18:H18").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Populate the newly inserted row 18 with data
$ws.Range("A18").Value = "1.1/2.0"
$ws.Range("B18").Value = "Yale"
$ws.Range("D18").Value = 5368
$ws.Range("E18").Value = "low"
$ws.Range("F18").Value = 5
$ws.Range("H18").Value = "Implement a CAS plugin for the Authentication Service."

# Row heights: the new data row is a bit taller; the blank spacer row (now 19)
# keeps its original short height.
$ws.Rows.Item(18).RowHeight = 16.5
$ws.Rows.Item(19).RowHeight = 13.5

$ws.Range("J8").Select() | Out-Null
